$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 1: "99.97" -> "0M", then insert 12 new rows right after it:
#     0M, 0M, 100, 0.00003, 0.00007, 0.00004, 0.00001, 0.00004, 0.00004,
#     0.00005, 0.00433, 100.0
$t.Cell(1,1).Range.Text = "0M"
$idx = 1
$vals1 = @("0M","0M","100","0.00003","0.00007","0.00004","0.00001","0.00004","0.00004","0.00005","0.00433","100.0")
foreach ($v in $vals1) {
    $idx = $idx + 1
    $newRow = $t.Rows.Add($t.Rows.Item($idx))
    $newRow.Cells.Item(1).Range.Text = $v
}

# --- Row 15 (was row 3 "15", now shifted down by the 12 rows above):
#     "15" -> "0.00000", then insert 9 new rows right after it:
#     0.00000 x8, then 0.0
$t.Cell(15,1).Range.Text = "0.00000"
$idx = 15
$vals2 = @("0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.0")
foreach ($v in $vals2) {
    $idx = $idx + 1
    $newRow = $t.Rows.Add($t.Rows.Item($idx))
    $newRow.Cells.Item(1).Range.Text = $v
}

# --- Row 45 (was row 24, the tab-separated "100 ... 100.0" row, now
#     shifted down by 12 + 9 = 21): collapse to the single value "99.97".
$t.Cell(45,1).Range.Text = "99.97"

# --- Row 46 (was row 25, the empty row): delete it entirely.
$t.Rows.Item(46).Delete()

# --- Row 46 (was row 26 "0[TAB]0...0", now shifted up by 1 because of the
#     deletion above): keep only the first "0", then append a new row "15".
$t.Cell(46,1).Range.Text = "0"
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "15"

Write-Output ("Final row count: " + $t.Rows.Count)
